$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$pairs = @(
    @(47, 48),
    @(55, 56),
    @(71, 72),
    @(109, 110),
    @(133, 134),
    @(149, 150),
    @(213, 214),
    @(229, 231),
    @(232, 233),
    @(245, 246),
    @(248, 249),
    @(263, 265),
    @(271, 272)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $range1 = $ws.Range("B$r1`:AB$r1")
    $range2 = $ws.Range("B$r2`:AB$r2")
    $v1 = $range1.Value2
    $v2 = $range2.Value2
    $range1.Value2 = $v2
    $range2.Value2 = $v1
}

Write-Host "Swapped $($pairs.Count) row pairs"
